$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4982.684
$ws.Range("I15").Value = 4982.684
$ws.Range("K15").Value = 14948.052
$ws.Range("M15").Value = -14779.052
$ws.Range("H17").Value = 2507.3157
$ws.Range("J17").Value = 2507.3157
$ws.Range("L17").Value = 7521.9471
$ws.Range("N17").Value = -7857.9471
$ws.Range("H116").Value = 3532.6667
$ws.Range("J116").Value = 3799.75
$ws.Range("L116").Value = 3799.75
$ws.Range("N116").Value = -10683.75
$ws.Range("H118").Value = 315.6154
$ws.Range("I118").Value = 300.25
$ws.Range("K118").Value = 900.75
$ws.Range("M118").Value = 756.25
$ws.Range("H121").Value = 5999
$ws.Range("J121").Value = 5999
$ws.Range("L121").Value = 17997
$ws.Range("N121").Value = -21491
$ws.Range("H131").Value = 4380.6875
$ws.Range("I131").Value = 3339.4
$ws.Range("K131").Value = 10018.2
$ws.Range("M131").Value = -4978.200000000001
$ws.Range("H138").Value = 3789.913
$ws.Range("I138").Value = 2197.1667
$ws.Range("J138").Value = 4352.0586
$ws.Range("K138").Value = 6591.500100000001
$ws.Range("L138").Value = 13056.1758
$ws.Range("M138").Value = -1451.500100000001
$ws.Range("N138").Value = -23336.1758

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3924.75
$ws.Range("I26").Value = 3924.75
$ws.Range("K26").Value = 3924.75
$ws.Range("M26").Value = -3594.75
$ws.Range("H32").Value = 9871.857
$ws.Range("I32").Value = 9422.571
$ws.Range("J32").Value = 11219.714
$ws.Range("K32").Value = 9422.571
$ws.Range("L32").Value = 11219.714
$ws.Range("M32").Value = -9135.571
$ws.Range("N32").Value = -11793.714
$ws.Range("H61").Value = 4953.1333
$ws.Range("I61").Value = 5126.024
$ws.Range("K61").Value = 5126.024
$ws.Range("M61").Value = -4914.024
$ws.Range("H74").Value = 3631.8948
$ws.Range("I74").Value = 1812.875
$ws.Range("J74").Value = 13333.333
$ws.Range("K74").Value = 1812.875
$ws.Range("L74").Value = 13333.333
$ws.Range("M74").Value = -938.875
$ws.Range("N74").Value = -15081.333
$ws.Range("H77").Value = 3631.8948
$ws.Range("I77").Value = 1812.875
$ws.Range("J77").Value = 13333.333
$ws.Range("K77").Value = 9064.375
$ws.Range("L77").Value = 66666.66500000001
$ws.Range("M77").Value = -4696.375
$ws.Range("N77").Value = -75402.66500000001
$ws.Range("H102").Value = 3312.3333
$ws.Range("I102").Value = 2115.8572
$ws.Range("K102").Value = 2115.8572
$ws.Range("M102").Value = -493.8571999999999
$ws.Range("H122").Value = 2147
$ws.Range("I122").Value = 2147
$ws.Range("K122").Value = 6441
$ws.Range("M122").Value = -3991
$ws.Range("H136").Value = 4953.1333
$ws.Range("I136").Value = 5126.024
$ws.Range("K136").Value = 15378.072
$ws.Range("M136").Value = -12828.072

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 15500.667
$ws.Range("I37").Value = 6500
$ws.Range("J37").Value = 20001
$ws.Range("K37").Value = 6500
$ws.Range("L37").Value = 20001
$ws.Range("M37").Value = -6363
$ws.Range("N37").Value = -20275
$ws.Range("H99").Value = 4397.885
$ws.Range("I99").Value = 3252.125
$ws.Range("K99").Value = 3252.125
$ws.Range("M99").Value = -1754.125
$ws.Range("H134").Value = 4041.6453
$ws.Range("J134").Value = 4490.8
$ws.Range("L134").Value = 13472.4
$ws.Range("N134").Value = -18542.4

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4428.357
$ws.Range("I16").Value = 3639.9
$ws.Range("J16").Value = 6399.5
$ws.Range("K16").Value = 3639.9
$ws.Range("L16").Value = 6399.5
$ws.Range("M16").Value = -3352.9
$ws.Range("N16").Value = -6973.5
$ws.Range("H58").Value = 11525.392
$ws.Range("I58").Value = 8095.857
$ws.Range("J58").Value = 13025.8125
$ws.Range("K58").Value = 8095.857
$ws.Range("L58").Value = 13025.8125
$ws.Range("M58").Value = -7892.857
$ws.Range("N58").Value = -13431.8125
$ws.Range("H86").Value = 8901.357
$ws.Range("I86").Value = 8102.5
$ws.Range("K86").Value = 8102.5
$ws.Range("M86").Value = -6979.5
$ws.Range("H89").Value = 8901.357
$ws.Range("I89").Value = 8102.5
$ws.Range("K89").Value = 40512.5
$ws.Range("M89").Value = -34896.5
$ws.Range("H105").Value = 3733.1667
$ws.Range("J105").Value = 3500
$ws.Range("L105").Value = 3500
$ws.Range("N105").Value = -6994
$ws.Range("H113").Value = 4428.357
$ws.Range("I113").Value = 3639.9
$ws.Range("J113").Value = 6399.5
$ws.Range("K113").Value = 3639.9
$ws.Range("L113").Value = 6399.5
$ws.Range("M113").Value = -1469.9
$ws.Range("N113").Value = -10739.5
$ws.Range("H136").Value = 11525.392
$ws.Range("I136").Value = 8095.857
$ws.Range("J136").Value = 13025.8125
$ws.Range("K136").Value = 24287.571
$ws.Range("L136").Value = 39077.4375
$ws.Range("M136").Value = -21737.571
$ws.Range("N136").Value = -44177.4375

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 25200062
$ws.Range("I11").Value = 25200062
$ws.Range("K11").Value = 75600186
$ws.Range("M11").Value = -75600046
$ws.Range("H12").Value = 210.75
$ws.Range("J12").Value = 325.4
$ws.Range("L12").Value = 976.1999999999999
$ws.Range("N12").Value = -1322.2
$ws.Range("H34").Value = 167044.17
$ws.Range("J34").Value = 200400
$ws.Range("L34").Value = 601200
$ws.Range("N34").Value = -601368
$ws.Range("H39").Value = 1601.3636
$ws.Range("J39").Value = 1777.7778
$ws.Range("L39").Value = 5333.3334
$ws.Range("N39").Value = -5921.3334
$ws.Range("H44").Value = 1734056.1
$ws.Range("I44").Value = 5200230
$ws.Range("J44").Value = 969.3
$ws.Range("K44").Value = 15600690
$ws.Range("L44").Value = 2907.9
$ws.Range("M44").Value = -15600292
$ws.Range("N44").Value = -3703.9
$ws.Range("H50").Value = 817.8182
$ws.Range("I50").Value = 947.1429000000001
$ws.Range("J50").Value = 591.5
$ws.Range("K50").Value = 2841.4287
$ws.Range("L50").Value = 1774.5
$ws.Range("M50").Value = -2360.4287
$ws.Range("N50").Value = -2736.5
$ws.Range("H53").Value = 817.8182
$ws.Range("I53").Value = 947.1429000000001
$ws.Range("J53").Value = 591.5
$ws.Range("K53").Value = 2841.4287
$ws.Range("L53").Value = 1774.5
$ws.Range("M53").Value = -2360.4287
$ws.Range("N53").Value = -2736.5
$ws.Range("H54").Value = 62675
$ws.Range("I54").Value = 233.33333
$ws.Range("J54").Value = 250000
$ws.Range("K54").Value = 699.99999
$ws.Range("L54").Value = 750000
$ws.Range("M54").Value = -140.99999
$ws.Range("N54").Value = -751118
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").Value = $null

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5409.4614
$ws.Range("I102").Value = 3923.2222
$ws.Range("K102").Value = 3923.2222
$ws.Range("M102").Value = -2301.2222
$ws.Range("H132").Value = 5267.913
$ws.Range("I132").Value = 4938.8945
$ws.Range("K132").Value = 14816.6835
$ws.Range("M132").Value = -12286.6835

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3382.6667
$ws.Range("I22").Value = 3549
$ws.Range("J22").Value = 3299.5
$ws.Range("K22").Value = 3549
$ws.Range("L22").Value = 3299.5
$ws.Range("M22").Value = -3254
$ws.Range("N22").Value = -3889.5
$ws.Range("H27").Value = 3382.6667
$ws.Range("I27").Value = 3549
$ws.Range("J27").Value = 3299.5
$ws.Range("K27").Value = 3549
$ws.Range("L27").Value = 3299.5
$ws.Range("M27").Value = -3442
$ws.Range("N27").Value = -3513.5
$ws.Range("H40").Value = 3821.111
$ws.Range("I40").Value = 3741.5715
$ws.Range("J40").Value = 4099.5
$ws.Range("K40").Value = 3741.5715
$ws.Range("L40").Value = 4099.5
$ws.Range("M40").Value = -3605.5715
$ws.Range("N40").Value = -4371.5
$ws.Range("H68").Value = 3291.35
$ws.Range("J68").Value = 5804
$ws.Range("L68").Value = 5804
$ws.Range("N68").Value = -7302
$ws.Range("H71").Value = 3291.35
$ws.Range("J71").Value = 5804
$ws.Range("L71").Value = 29020
$ws.Range("N71").Value = -36508
$ws.Range("H122").Value = 2621.3333
$ws.Range("J122").Value = 2699
$ws.Range("L122").Value = 8097
$ws.Range("N122").Value = -12997
$ws.Range("H132").Value = 7198.8203
$ws.Range("I132").Value = 8151.9
$ws.Range("K132").Value = 24455.7
$ws.Range("M132").Value = -21925.7

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 646.4400000000001
$ws.Range("I113").Value = 500.46667
$ws.Range("J113").Value = 865.4
$ws.Range("K113").Value = 1501.40001
$ws.Range("L113").Value = 2596.2
$ws.Range("M113").Value = 668.5999899999999
$ws.Range("N113").Value = -6936.2
